# Applies the "Accomplish PDF download - first try" edit:
# - Update the PDF filename for row 2
# - Mark rows 3 and 4 as "Sent" = Yes
# - Replace Micaela's row (row 5) with Matias's data
# - Add a new row 6 for Pablo

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update PDF file name
$ws.Range("C2").Value = "Pdf_de_prueba.pdf"

# Row 3 & 4: mark as sent
$ws.Range("D3").Value = "Yes"
$ws.Range("D4").Value = "Yes"

# Row 5: replace Micaela with Matias
$ws.Range("A5").Value = "Matias"
$ws.Range("B5").Value = "matiasmalleville@gmail.com"
$ws.Range("C5").Value = "Matias Malleville.pdf"
$ws.Range("D5").Value = "Yes"

# Row 6: new recipient Pablo
$ws.Range("A6").Value = "Pablo"
$ws.Range("B6").Value = "roig@lacaja.com.ar"
$ws.Range("C6").Value = "Pablo Roig.pdf"
$ws.Range("D6").Value = "Yes"
